$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7899
$ws.Range("L2").Value = 732
$ws.Range("L3").Value = 733
$ws.Range("L4").Value = 185
$ws.Range("D6").Value = 11881
$ws.Range("K6").Value = 9122
$ws.Range("L6").Value = 789
$ws.Range("D7").Value = 28176
$ws.Range("K7").Value = 27521
$ws.Range("L7").Value = 2492

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 41
$ws.Range("L3").Value = 49
$ws.Range("L6").Value = 51
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 22
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L2").Value = 14
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 11
$ws.Range("K6").Value = 113
$ws.Range("K7").Value = 464
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L4").Value = 9
$ws.Range("L6").Value = 21
$ws.Range("L8").Value = 155
$ws.Range("L11").Value = 39
$ws.Range("L17").Value = 5
$ws.Range("L19").Value = 82
$ws.Range("L29").Value = 126
$ws.Range("L32").Value = 5
$ws.Range("L33").Value = 104
$ws.Range("L36").Value = 46
$ws.Range("L37").Value = 80
$ws.Range("L41").Value = 12
$ws.Range("L42").Value = 84
$ws.Range("L43").Value = 18
$ws.Range("L49").Value = 15
$ws.Range("L50").Value = 20
$ws.Range("K52").Value = 700
$ws.Range("L52").Value = 47
$ws.Range("L54").Value = 51
$ws.Range("L59").Value = 2
$ws.Range("D63").Value = 365
$ws.Range("L63").Value = 12
$ws.Range("L67").Value = 87
$ws.Range("L68").Value = 6
$ws.Range("L73").Value = 16
$ws.Range("L78").Value = 32
$ws.Range("L84").Value = 21
$ws.Range("L85").Value = 121
$ws.Range("L86").Value = 17
$ws.Range("L88").Value = 40
$ws.Range("L93").Value = 15
$ws.Range("L95").Value = 35
$ws.Range("L97").Value = 30
$ws.Range("K99").Value = 464
$ws.Range("L99").Value = 41
$ws.Range("D101").Value = 28176
$ws.Range("K101").Value = 27521
$ws.Range("L101").Value = 2492

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 22
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 45
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 82

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 84

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 21
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L2").Value = 1
$ws.Range("L7").Value = 2

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L4").ClearContents()
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 5

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 6

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 56
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 191
$ws.Range("L2").Value = 17
$ws.Range("L3").Value = 13
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 9
